# "length matched DDR_A byte 1"
# Re-run of trace length matching for the DDR_A "Byte Lane 1" net group
# (rows 4-14 of the DDR_A sheet): Track Length (E), Total Length (G),
# Track Delay (H) and Total Delay (L) all move together as the traces
# were re-routed to match length; Via Length (F), Via Delay (I), Package
# Delay (J) and Extra Delay (K) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DDR_A")

# row -> (TrackLength E, TotalLength G, TrackDelay H, TotalDelay L)
$rows = @(
    @{ Row = 4;  E = 11.44; G = 13.84; H = 67.03;  L = 151.63 },
    @{ Row = 5;  E = 11.71; G = 14.11; H = 68.6;    L = 151.8  },
    @{ Row = 6;  E = 17.58; G = 17.58; H = 102.83; L = 151.63 },
    @{ Row = 7;  E = 17;    G = 17;    H = 99.31;  L = 151.71 },
    @{ Row = 8;  E = 15.97; G = 15.97; H = 93.39;  L = 151.79 },
    @{ Row = 9;  E = 15.9;  G = 15.9;  H = 92.97;  L = 151.57 },
    @{ Row = 10; E = 18.25; G = 18.25; H = 106.62; L = 151.62 },
    @{ Row = 11; E = 18.03; G = 18.03; H = 105.4;  L = 151.6  },
    @{ Row = 12; E = 17.36; G = 17.36; H = 101.44; L = 151.54 },
    @{ Row = 13; E = 17.86; G = 17.86; H = 104.38; L = 151.58 },
    @{ Row = 14; E = 19.06; G = 19.06; H = 111.3;  L = 151.6  }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("E$n").Value = $r.E
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    $ws.Range("L$n").Value = $r.L
}

# The edit was made while DDR_A was the active sheet; the cursor ended up
# on A36:B36 once the byte-lane edits were done.
$ws.Activate() | Out-Null
$ws.Range("A36:B36").Select() | Out-Null
